{"js": "// Update the date heading and every two-digit multiplication answer cell in the\n// table. Every old value below is unique within the document, so an exact,\n// case-sensitive search reliably targets the single run that needs updating.\nconst replacements = [\n  [\"2023-08-11 Friday\", \"2023-08-12 Saturday\"],\n  [\"89\u00d721=1869\", \"17\u00d755=935\"],\n  [\"78\u00d775=5850\", \"12\u00d737=444\"],\n  [\"60\u00d726=1560\", \"97\u00d713=1261\"],\n  [\"91\u00d751=4641\", \"98\u00d744=4312\"],\n  [\"39\u00d776=2964\", \"28\u00d738=1064\"],\n  [\"94\u00d749=4606\", \"20\u00d742=840\"],\n  [\"82\u00d749=4018\", \"96\u00d718=1728\"],\n  [\"76\u00d742=3192\", \"20\u00d787=1740\"],\n  [\"30\u00d763=1890\", \"74\u00d777=5698\"],\n  [\"57\u00d750=2850\", \"38\u00d737=1406\"],\n  [\"88\u00d760=5280\", \"21\u00d781=1701\"],\n  [\"57\u00d719=1083\", \"74\u00d751=3774\"],\n  [\"74\u00d791=6734\", \"51\u00d750=2550\"],\n  [\"65\u00d727=1755\", \"43\u00d727=1161\"],\n  [\"58\u00d745=2610\", \"85\u00d770=5950\"],\n  [\"78\u00d781=6318\", \"87\u00d753=4611\"],\n  [\"66\u00d716=1056\", \"93\u00d738=3534\"],\n  [\"46\u00d737=1702\", \"42\u00d720=840\"],\n  [\"12\u00d777=924\", \"57\u00d770=3990\"],\n  [\"62\u00d799=6138\", \"15\u00d776=1140\"],\n  [\"60\u00d757=3420\", \"76\u00d739=2964\"],\n  [\"11\u00d760=660\", \"16\u00d762=992\"],\n  [\"41\u00d799=4059\", \"23\u00d735=805\"],\n  [\"40\u00d794=3760\", \"50\u00d780=4000\"],\n  [\"83\u00d750=4150\", \"41\u00d711=451\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old value with its corresponding new value.\n# All source strings in this document are unique, so an exact,\n# case-sensitive match targets exactly one run each.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2023-08-11 Friday\", \"2023-08-12 Saturday\")\n    ,@(\"89\u00d721=1869\", \"17\u00d755=935\")\n    ,@(\"78\u00d775=5850\", \"12\u00d737=444\")\n    ,@(\"60\u00d726=1560\", \"97\u00d713=1261\")\n    ,@(\"91\u00d751=4641\", \"98\u00d744=4312\")\n    ,@(\"39\u00d776=2964\", \"28\u00d738=1064\")\n    ,@(\"94\u00d749=4606\", \"20\u00d742=840\")\n    ,@(\"82\u00d749=4018\", \"96\u00d718=1728\")\n    ,@(\"76\u00d742=3192\", \"20\u00d787=1740\")\n    ,@(\"30\u00d763=1890\", \"74\u00d777=5698\")\n    ,@(\"57\u00d750=2850\", \"38\u00d737=1406\")\n    ,@(\"88\u00d760=5280\", \"21\u00d781=1701\")\n    ,@(\"57\u00d719=1083\", \"74\u00d751=3774\")\n    ,@(\"74\u00d791=6734\", \"51\u00d750=2550\")\n    ,@(\"65\u00d727=1755\", \"43\u00d727=1161\")\n    ,@(\"58\u00d745=2610\", \"85\u00d770=5950\")\n    ,@(\"78\u00d781=6318\", \"87\u00d753=4611\")\n    ,@(\"66\u00d716=1056\", \"93\u00d738=3534\")\n    ,@(\"46\u00d737=1702\", \"42\u00d720=840\")\n    ,@(\"12\u00d777=924\", \"57\u00d770=3990\")\n    ,@(\"62\u00d799=6138\", \"15\u00d776=1140\")\n    ,@(\"60\u00d757=3420\", \"76\u00d739=2964\")\n    ,@(\"11\u00d760=660\", \"16\u00d762=992\")\n    ,@(\"41\u00d799=4059\", \"23\u00d735=805\")\n    ,@(\"40\u00d794=3760\", \"50\u00d780=4000\")\n    ,@(\"83\u00d750=4150\", \"41\u00d711=451\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap:=wdFindContinue, Format, ReplaceWith,\n    #          Replace:=wdReplaceOne)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n}\n"}
